$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column B (pushes existing B:C to C:D) ---
$ws.Columns.Item(2).Insert()

# --- Column width for the new column B (match column A) ---
$ws.Columns.Item(2).ColumnWidth = 75.81640625

# --- Row 2 height (tall row for wrapped multi-line query text) ---
$ws.Rows.Item(2).RowHeight = 101.5

# --- Header row ---
$ws.Range("B1").Value = "StatQuery"

# --- Row 2 query cells (cells already inherit the wrap-text style from the insert) ---
$statsQuery = @'
MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report)OPTIONAL MATCH (s)<-[*]-(f:file)  WITH  c AS c, t ,a, s , f WHERE f.file_type IN ['Aligned DNA reads file','Aligned RNA reads file','Index file','Variants file'] RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trials
'@

$indexFileQuery = @'
MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report)OPTIONAL MATCH (s)<-[*]-(f:file)  WITH  c AS c, t ,a, s WHERE f.file_type IN ['Index file']  RETURN DISTINCT coalesce(c.case_id,'') AS `Case ID` , coalesce(t.clinical_trial_designation ,'')as `Trial Code` , coalesce(a.arm_id,'') As `Arm` , coalesce(a.arm_drug,'') As `Arm Treatment` , coalesce(c.disease,'') As Diagnosis , coalesce(c.gender,'') As Gender , coalesce(c.race,'') As Race , coalesce(c.ethnicity,'') As Ethnicity
'@

$ws.Range("A2").Value = $indexFileQuery
$ws.Range("B2").Value = $statsQuery

# --- Selection / scroll position ---
$ws.Activate()
$ws.Range("A7").Select()
